# Fix column name in quantitation example file:
#   D1 header: "standard_ ng" (typo, stray space) -> "standard_ng"
#   E1 header: "sample_ul" stays the same value (only its position in the
#              shared-string table shifts as a side effect of the rename)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "standard_ng"
$ws.Range("E1").Value = "sample_ul"
